$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "42.033.37"
$cell.ClearFormats()
$ws.Range("E2").Value = "  -0.71%  "
$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "2.247.54"
$cell.ClearFormats()
$ws.Range("E3").Value = "  -1.28%  "
$ws.Range("E4").Value = "  -0.15%  "
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "306.40"
$cell.ClearFormats()
$ws.Range("E5").Value = "  -0.14%  "
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "96.30"
$cell.ClearFormats()
$ws.Range("E6").Value = "  -1.35%  "
$ws.Range("E7").Value = "  -1.29%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("E9").Value = "  -1.46%  "
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "34.68"
$cell.ClearFormats()
$ws.Range("E10").Value = "  -2.64%  "
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "0.0809"
$cell.ClearFormats()
$ws.Range("E11").Value = "  +1.57%  "
$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "0.113"
$cell.ClearFormats()
$ws.Range("E12").Value = "  +0.42%  "
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "6.77"
$cell.ClearFormats()
$ws.Range("E13").Value = "  +1.37%  "
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "2.597.27"
$cell.ClearFormats()
$ws.Range("E14").Value = "  -1.28%  "
$ws.Range("E15").Value = "  +0.07%  "
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "2.237.18"
$cell.ClearFormats()
$ws.Range("E16").Value = "  -1.53%  "
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "0.777"
$cell.ClearFormats()
$ws.Range("E17").Value = "  -2.29%  "
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "41.896.30"
$cell.ClearFormats()
$ws.Range("E18").Value = "  -0.85%  "
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "12.17"
$cell.ClearFormats()
$ws.Range("E19").Value = "  -2.88%  "
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "0.0₃0901"
$cell.ClearFormats()
$ws.Range("E20").Value = "  -0.99%  "
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "5.91"
$cell.ClearFormats()
$ws.Range("E21").Value = "  -0.88%  "
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "67.14"
$cell.ClearFormats()
$ws.Range("E22").Value = "  -0.73%  "
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "235.25"
$cell.ClearFormats()
$ws.Range("E23").Value = "  -2.36%  "
$ws.Range("E24").Value = "  -1.33%  "
$ws.Range("E25").Value = "  -0.19%  "
$ws.Range("E26").Value = "  +0.00%  "
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "23.30"
$cell.ClearFormats()
$ws.Range("E27").Value = "  -2.11%  "
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "36.75"
$cell.ClearFormats()
$ws.Range("E28").Value = "  -2.17%  "
$ws.Range("B29").Value = "Cosmos"
$ws.Range("C29").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "9.48"
$cell.ClearFormats()
$ws.Range("E29").Value = "  -0.12%  "
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "2.11"
$cell.ClearFormats()
$ws.Range("E30").Value = "  +0.81%  "
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "164.89"
$cell.ClearFormats()
$ws.Range("E31").Value = "  +3.27%  "
$ws.Range("E32").Value = "  -0.02%  "
$ws.Range("E33").Value = "  -1.48%  "
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "3.08"
$cell.ClearFormats()
$ws.Range("E34").Value = "  -1.40%  "
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "17.58"
$cell.ClearFormats()
$ws.Range("E35").Value = "  +3.55%  "
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "0.0720"
$cell.ClearFormats()
$ws.Range("E36").Value = "  -2.76%  "
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "2.38"
$cell.ClearFormats()
$ws.Range("E37").Value = "  -0.51%  "
$ws.Range("E38").Value = "  -0.42%  "
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "0.102"
$cell.ClearFormats()
$ws.Range("E39").Value = "  -3.20%  "
$ws.Range("E40").Value = "  -2.97%  "
$ws.Range("E41").Value = "  -1.10%  "
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "1.940.44"
$cell.ClearFormats()
$ws.Range("E42").Value = "  -3.04%  "
$ws.Range("E43").Value = "  -1.72%  "
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "2.17"
$cell.ClearFormats()
$ws.Range("E44").Value = "  -9.80%  "
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "18.44"
$cell.ClearFormats()
$ws.Range("E45").Value = "  -3.04%  "
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "2.90"
$cell.ClearFormats()
$ws.Range("E46").Value = "  -3.15%  "
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "9.65"
$cell.ClearFormats()
$ws.Range("E47").Value = "  -3.36%  "
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "53.52"
$cell.ClearFormats()
$ws.Range("E48").Value = "  +0.99%  "
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "2.469.51"
$cell.ClearFormats()
$ws.Range("E49").Value = "  -1.16%  "
$ws.Range("E50").Value = "  -1.37%  "
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "91.02"
$cell.ClearFormats()
$ws.Range("E51").Value = "  -0.92%  "

Write-Output "Applied 89 cell updates"
